$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")
$ws.Columns.Item(7).Delete()
